$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental row: set the "Value" cell (B7) to the literal text "true"
# (a plain string, not a native Excel boolean). Assigning the string
# directly via .Value would be auto-coerced to a TRUE boolean, so instead
# compute it with a text formula and then flatten the formula down to a
# static value via copy / paste-values (keeps the original cell style).
$ws.Range("B7").Formula = '=TEXT("true","@")'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Bump the generation Date value to the new timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
